# Update "想去人数" (interest count) values in the F column on both the
# "展览" sheet and the "全部类型" aggregate sheet, as produced by the
# gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1156
$ws1.Range("F3").Value = 1969
$ws1.Range("F4").Value = 623
$ws1.Range("F5").Value = 1276
$ws1.Range("F12").Value = 858
$ws1.Range("F17").Value = 349
$ws1.Range("F22").Value = 205
$ws1.Range("F25").Value = 375

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1156
$ws4.Range("F4").Value = 1969
$ws4.Range("F5").Value = 623
$ws4.Range("F6").Value = 1276
$ws4.Range("F14").Value = 858
$ws4.Range("F22").Value = 349
$ws4.Range("F29").Value = 205
$ws4.Range("F32").Value = 375
